{"js": "// 1. Insert the new paragraph right after the \"Nedan presenteras ...\" intro\n//    paragraph near the top of the document (it is currently duplicated at\n//    the very end of the document body, which is removed in step 2).\nconst body = context.document.body;\nlet paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst introParagraph = paragraphs.items[2];\nintroParagraph.insertParagraph(\n  \"Vi f\u00f6rv\u00e4ntar oss att ni \u00e5terkommer med ett skriftligt svar p\u00e5 v\u00e5rt klagom\u00e5l och \u00e4ven beskriver vilka korrigerande \u00e5tg\u00e4rder ni satt in f\u00f6r att r\u00e4tta till identifierade brister i er efterlevnad av den svenska FSC standarden.\",\n  \"After\"\n);\nawait context.sync();\n\n// 2. Remove the trailing two empty paragraphs and the old copy of that same\n//    paragraph that used to sit at the end of the document body. After the\n//    insertion above, these are still the last three paragraphs, so delete\n//    them by repeatedly removing the paragraph now three-from-the-end.\nparagraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nlet count = paragraphs.items.length;\nconst deleteFromIndex = count - 3;\n\nfor (let i = 0; i < 3; i++) {\n  paragraphs = body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n  paragraphs.items[deleteFromIndex].delete();\n  await context.sync();\n}\n\n// 3. Update the date shown in the first-page header from 2023-11-13 to\n//    2023-11-14.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\nconst firstPageHeader = sections.items[0].getHeader(\"FirstPage\");\n\nconst found = firstPageHeader.search(\"2023-11-13\");\nfound.load(\"items\");\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].insertText(\"2023-11-14\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Insert the new paragraph right after the \"Nedan presenteras ...\" paragraph\n#    near the top of the document (it is currently duplicated at the very\n#    end of the document, which we remove in step 2).\n$introPar = $d.Paragraphs.Item(3)\n$introPar.Range.InsertParagraphAfter()\n$newPar = $d.Paragraphs.Item(4)\n$newPar.Range.Text = \"Vi f\u00f6rv\u00e4ntar oss att ni \u00e5terkommer med ett skriftligt svar p\u00e5 v\u00e5rt klagom\u00e5l och \u00e4ven beskriver vilka korrigerande \u00e5tg\u00e4rder ni satt in f\u00f6r att r\u00e4tta till identifierade brister i er efterlevnad av den svenska FSC standarden.\"\n\n# 2. Remove the trailing two empty paragraphs and the old copy of that\n#    same paragraph at the end of the document body.\n$count = $d.Paragraphs.Count\n$startPar = $d.Paragraphs.Item($count - 2)\n$endPar = $d.Paragraphs.Item($count)\n$rng = $d.Range($startPar.Range.Start, $endPar.Range.End)\n$rng.Delete()\n\n# 3. Update the date shown in the first-page header from 2023-11-13 to\n#    2023-11-14.\n$section = $d.Sections.First\n$firstPageHeader = $section.Headers.Item(2)\n$firstPageHeader.Range.Find.Execute(\"2023-11-13\", $false, $false, $false, $false, $false, $true, 1, $false, \"2023-11-14\", 2)\n"}
